# The Chelyabinsk light-curve sheet had a stray blank row at row 181
# (no row 181 existed at all in the XML, so everything from 182-341 was
# effectively shifted down by one against the intended row numbering).
# This was throwing off the integrated-intensity calculations, so we
# remove that gap by deleting row 181, which shifts rows 182-341 up to
# occupy 181-340.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(181).Delete()

# Restore the view: scroll so row 163 is the top-left visible row, and
# leave the active selection on F182 (matching the author's viewport at
# the time of the edit).
$win = $excel.ActiveWindow
$ws.Range("F182").Select()
$win.ScrollRow = 163
$win.ScrollColumn = 1
